$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3098.8
$ws.Range("J43").Value = 3165.3333
$ws.Range("L43").Value = 3165.3333
$ws.Range("N43").Value = -3303.3333

$ws.Range("H53").Value = 423.3125
$ws.Range("I53").Value = 529.25
$ws.Range("J53").Value = 317.375
$ws.Range("K53").Value = 529.25
$ws.Range("L53").Value = 317.375
$ws.Range("M53").Value = 107.75
$ws.Range("N53").Value = -1591.375

$ws.Range("H92").Value = 4699.5835
$ws.Range("I92").Value = 4600.263
$ws.Range("J92").Value = 5077
$ws.Range("K92").Value = 4600.263
$ws.Range("L92").Value = 5077
$ws.Range("M92").Value = -3352.263
$ws.Range("N92").Value = -7573

$ws.Range("H132").Value = 1462.1538
$ws.Range("I132").Value = 1194.3226
$ws.Range("K132").Value = 3582.9678
$ws.Range("M132").Value = -1052.9678

$ws.Range("H135").Value = 1282.2727
$ws.Range("I135").Value = 1321.6072
$ws.Range("J135").Value = 1062
$ws.Range("K135").Value = 11894.4648
$ws.Range("L135").Value = 9558
$ws.Range("M135").Value = -9359.4648
$ws.Range("N135").Value = -14628

$ws.Range("H136").Value = 120064.086
$ws.Range("I136").Value = 80000
$ws.Range("J136").Value = 133418.78
$ws.Range("K136").Value = 80000
$ws.Range("L136").Value = 133418.78
$ws.Range("M136").Value = -74900
$ws.Range("N136").Value = -143618.78

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8980.924999999999
$ws.Range("I32").Value = 8636.579
$ws.Range("K32").Value = 8636.579
$ws.Range("M32").Value = -8349.579

$ws.Range("H45").Value = 1512.8334
$ws.Range("I45").Value = 1515.4
$ws.Range("K45").Value = 1515.4
$ws.Range("M45").Value = -1138.4

$ws.Range("H97").Value = 1877.6522
$ws.Range("I97").Value = 1335.75
$ws.Range("K97").Value = 1335.75
$ws.Range("M97").Value = -839.75

$ws.Range("H132").Value = 6908.933
$ws.Range("I132").Value = 5505.4
$ws.Range("K132").Value = 16516.2
$ws.Range("M132").Value = -13986.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 41535.777
$ws.Range("I75").Value = 40477.75
$ws.Range("J75").Value = 50000
$ws.Range("K75").Value = 40477.75
$ws.Range("L75").Value = 50000
$ws.Range("M75").Value = -39541.75
$ws.Range("N75").Value = -51872

$ws.Range("H78").Value = 41535.777
$ws.Range("I78").Value = 40477.75
$ws.Range("J78").Value = 50000
$ws.Range("K78").Value = 121433.25
$ws.Range("L78").Value = 150000
$ws.Range("M78").Value = -116753.25
$ws.Range("N78").Value = -159360

$ws.Range("H105").Value = 1749.375
$ws.Range("J105").Value = 939.2
$ws.Range("L105").Value = 939.2
$ws.Range("N105").Value = -4433.2

$ws.Range("H107").Value = 2353.44
$ws.Range("I107").Value = 1922.4
$ws.Range("K107").Value = 1922.4
$ws.Range("M107").Value = -2.400000000000091

$ws.Range("H134").Value = 2871.1428
$ws.Range("I134").Value = 2364
$ws.Range("J134").Value = 4494
$ws.Range("K134").Value = 7092
$ws.Range("L134").Value = 13482
$ws.Range("M134").Value = -4557
$ws.Range("N134").Value = -18552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 46332
$ws.Range("J74").Value = 46332
$ws.Range("L74").Value = 46332
$ws.Range("N74").Value = -48080

$ws.Range("H77").Value = 46332
$ws.Range("J77").Value = 46332
$ws.Range("L77").Value = 138996
$ws.Range("N77").Value = -147732

$ws.Range("H132").Value = 3295.9375
$ws.Range("I132").Value = 2378.825
$ws.Range("K132").Value = 7136.474999999999
$ws.Range("M132").Value = -4606.474999999999

$ws.Range("H134").Value = 5072.758
$ws.Range("I134").Value = 5406.357
$ws.Range("K134").Value = 16219.071
$ws.Range("M134").Value = -13684.071

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 470
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9590.4
$ws.Range("I70").Value = 8918.286
$ws.Range("K70").Value = 8918.286
$ws.Range("M70").Value = -8648.286

$ws.Range("H73").Value = 9590.4
$ws.Range("I73").Value = 8918.286
$ws.Range("K73").Value = 8918.286
$ws.Range("M73").Value = -7982.286

$ws.Range("H80").Value = 2127.383
$ws.Range("I80").Value = 2222.0908
$ws.Range("K80").Value = 2222.0908
$ws.Range("M80").Value = -1224.0908

$ws.Range("H83").Value = 2127.383
$ws.Range("I83").Value = 2222.0908
$ws.Range("K83").Value = 11110.454
$ws.Range("M83").Value = -6118.454

$ws.Range("H97").Value = 4273.25
$ws.Range("I97").Value = 4054.4783
$ws.Range("J97").Value = 5279.6
$ws.Range("K97").Value = 4054.4783
$ws.Range("L97").Value = 5279.6
$ws.Range("M97").Value = -3558.4783
$ws.Range("N97").Value = -6271.6

$ws.Range("H126").Value = 10394.223
$ws.Range("I126").Value = 12050
$ws.Range("J126").Value = 9069.6
$ws.Range("K126").Value = 36150
$ws.Range("L126").Value = 27208.8
$ws.Range("M126").Value = -33680
$ws.Range("N126").Value = -32148.8

$ws.Range("H132").Value = 74136.875
$ws.Range("I132").Value = 106890.3
$ws.Range("K132").Value = 320670.9
$ws.Range("M132").Value = -318140.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5284.868
$ws.Range("I132").Value = 4050
$ws.Range("J132").Value = 6567.231
$ws.Range("K132").Value = 12150
$ws.Range("L132").Value = 19701.693
$ws.Range("M132").Value = -9620
$ws.Range("N132").Value = -24761.693

$ws.Range("H136").Value = 3305.6667
$ws.Range("I136").Value = 3468.1924
$ws.Range("J136").Value = 2249.25
$ws.Range("K136").Value = 10404.5772
$ws.Range("L136").Value = 6747.75
$ws.Range("M136").Value = -7854.5772
$ws.Range("N136").Value = -11847.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4360.579
$ws.Range("I122").Value = 4190.7334
$ws.Range("K122").Value = 12572.2002
$ws.Range("M122").Value = -10122.2002

$ws.Range("H132").Value = 1319.6471
$ws.Range("I132").Value = 1188.6666
$ws.Range("J132").Value = 2302
$ws.Range("K132").Value = 3565.9998
$ws.Range("L132").Value = 6906
$ws.Range("M132").Value = -1035.9998
$ws.Range("N132").Value = -11966

Write-Host "Done applying Lich_Profits updates"
